{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst bangs = \"!\".repeat(48);\n\n// Replace the first paragraph's text with the bang string, keeping the\n// paragraph itself (and its formatting) intact.\nparagraphs.items[0].insertText(bangs, Word.InsertLocation.replace);\n\n// Remove the now-redundant second paragraph entirely.\nparagraphs.items[1].delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 48 exclamation marks - the new content for the first paragraph.\n$bangs = \"!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!!\"\n\n# Overwrite the first paragraph's content (leaving its own trailing\n# paragraph mark / formatting untouched) with the bangs text.\n$d.Paragraphs(1).Range.Text = $bangs\n\n# The second paragraph (\"I still wait for your code!\") is no longer\n# needed - remove it (and its paragraph mark) entirely.\n$d.Paragraphs(2).Range.Delete()\n"}
